$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 213
$ws.Range("A213").Formula = "=SUM(A212,1)"
$ws.Range("B213").Value = "preferredLang"
$ws.Range("C213").Value = "user preferred Language"
$ws.Range("D213").Value = "string"
$ws.Range("E213").Value = "{`"value`":`"Española`",`"code`":`"spa`"}"
$ws.Range("F213").Value = "eng"
$ws.Range("G213").Value = "TRUE"
$ws.Range("H213").Value = "superadmin"
$ws.Range("I213").Value = "now()"

# Row 214
$ws.Range("A214").Formula = "=SUM(A213,1)"
$ws.Range("B214").Value = "preferredLang"
$ws.Range("C214").Value = "user preferred Language"
$ws.Range("D214").Value = "string"
$ws.Range("E214").Value = "{`"value`":`"ಕನ್ನಡ`",`"code`":`"kan`"}"
$ws.Range("F214").Value = "eng"
$ws.Range("G214").Value = "TRUE"
$ws.Range("H214").Value = "superadmin"
$ws.Range("I214").Value = "now()"

# Row 215
$ws.Range("A215").Formula = "=SUM(A214,1)"
$ws.Range("B215").Value = "preferredLang"
$ws.Range("C215").Value = "user preferred Language"
$ws.Range("D215").Value = "string"
$ws.Range("E215").Value = "{`"value`":`"हिंदी`",`"code`":`"hin`"}"
$ws.Range("F215").Value = "eng"
$ws.Range("G215").Value = "TRUE"
$ws.Range("H215").Value = "superadmin"
$ws.Range("I215").Value = "now()"

# Row 216
$ws.Range("A216").Formula = "=SUM(A215,1)"
$ws.Range("B216").Value = "preferredLang"
$ws.Range("C216").Value = "user preferred Language"
$ws.Range("D216").Value = "string"
$ws.Range("E216").Value = "{`"value`":`"தமிழ்`",`"code`":`"tam`"}"
$ws.Range("F216").Value = "eng"
$ws.Range("G216").Value = "TRUE"
$ws.Range("H216").Value = "superadmin"
$ws.Range("I216").Value = "now()"

# Row 217
$ws.Range("A217").Formula = "=SUM(A216,1)"
$ws.Range("B217").Value = "preferredLang"
$ws.Range("C217").Value = "Langue préférée de l'utilisateur"
$ws.Range("D217").Value = "string"
$ws.Range("E217").Value = "{`"value`":`"Española`",`"code`":`"spa`"}"
$ws.Range("F217").Value = "fra"
$ws.Range("G217").Value = "TRUE"
$ws.Range("H217").Value = "superadmin"
$ws.Range("I217").Value = "now()"

# Row 218
$ws.Range("A218").Formula = "=SUM(A217,1)"
$ws.Range("B218").Value = "preferredLang"
$ws.Range("C218").Value = "Langue préférée de l'utilisateur"
$ws.Range("D218").Value = "string"
$ws.Range("E218").Value = "{`"value`":`"ಕನ್ನಡ`",`"code`":`"kan`"}"
$ws.Range("F218").Value = "fra"
$ws.Range("G218").Value = "TRUE"
$ws.Range("H218").Value = "superadmin"
$ws.Range("I218").Value = "now()"

# Row 219
$ws.Range("A219").Formula = "=SUM(A218,1)"
$ws.Range("B219").Value = "preferredLang"
$ws.Range("C219").Value = "Langue préférée de l'utilisateur"
$ws.Range("D219").Value = "string"
$ws.Range("E219").Value = "{`"value`":`"हिंदी`",`"code`":`"hin`"}"
$ws.Range("F219").Value = "fra"
$ws.Range("G219").Value = "TRUE"
$ws.Range("H219").Value = "superadmin"
$ws.Range("I219").Value = "now()"

# Row 220
$ws.Range("A220").Formula = "=SUM(A219,1)"
$ws.Range("B220").Value = "preferredLang"
$ws.Range("C220").Value = "Langue préférée de l'utilisateur"
$ws.Range("D220").Value = "string"
$ws.Range("E220").Value = "{`"value`":`"தமிழ்`",`"code`":`"tam`"}"
$ws.Range("F220").Value = "fra"
$ws.Range("G220").Value = "TRUE"
$ws.Range("H220").Value = "superadmin"
$ws.Range("I220").Value = "now()"

# Row 221
$ws.Range("A221").Formula = "=SUM(A220,1)"
$ws.Range("B221").Value = "preferredLang"
$ws.Range("C221").Value = "يفضل المستخدم اللغة"
$ws.Range("D221").Value = "string"
$ws.Range("E221").Value = "{`"value`":`"Española`",`"code`":`"spa`"}"
$ws.Range("F221").Value = "ara"
$ws.Range("G221").Value = "TRUE"
$ws.Range("H221").Value = "superadmin"
$ws.Range("I221").Value = "now()"

# Row 222
$ws.Range("A222").Formula = "=SUM(A221,1)"
$ws.Range("B222").Value = "preferredLang"
$ws.Range("C222").Value = "يفضل المستخدم اللغة"
$ws.Range("D222").Value = "string"
$ws.Range("E222").Value = "{`"value`":`"ಕನ್ನಡ`",`"code`":`"kan`"}"
$ws.Range("F222").Value = "ara"
$ws.Range("G222").Value = "TRUE"
$ws.Range("H222").Value = "superadmin"
$ws.Range("I222").Value = "now()"

# Row 223
$ws.Range("A223").Formula = "=SUM(A222,1)"
$ws.Range("B223").Value = "preferredLang"
$ws.Range("C223").Value = "يفضل المستخدم اللغة"
$ws.Range("D223").Value = "string"
$ws.Range("E223").Value = "{`"value`":`"हिंदी`",`"code`":`"hin`"}"
$ws.Range("F223").Value = "ara"
$ws.Range("G223").Value = "TRUE"
$ws.Range("H223").Value = "superadmin"
$ws.Range("I223").Value = "now()"

# Row 224
$ws.Range("A224").Formula = "=SUM(A223,1)"
$ws.Range("B224").Value = "preferredLang"
$ws.Range("C224").Value = "يفضل المستخدم اللغة"
$ws.Range("D224").Value = "string"
$ws.Range("E224").Value = "{`"value`":`"தமிழ்`",`"code`":`"tam`"}"
$ws.Range("F224").Value = "ara"
$ws.Range("G224").Value = "TRUE"
$ws.Range("H224").Value = "superadmin"
$ws.Range("I224").Value = "now()"

# Row 225
$ws.Range("A225").Formula = "=SUM(A224,1)"
$ws.Range("B225").Value = "preferredLang"
$ws.Range("C225").Value = "idioma preferido del usuario"
$ws.Range("D225").Value = "string"
$ws.Range("E225").Value = "{`"value`":`"ಕನ್ನಡ`",`"code`":`"kan`"}"
$ws.Range("F225").Value = "spa"
$ws.Range("G225").Value = "TRUE"
$ws.Range("H225").Value = "superadmin"
$ws.Range("I225").Value = "now()"

# Row 226
$ws.Range("A226").Formula = "=SUM(A225,1)"
$ws.Range("B226").Value = "preferredLang"
$ws.Range("C226").Value = "idioma preferido del usuario"
$ws.Range("D226").Value = "string"
$ws.Range("E226").Value = "{`"value`":`"हिंदी`",`"code`":`"hin`"}"
$ws.Range("F226").Value = "spa"
$ws.Range("G226").Value = "TRUE"
$ws.Range("H226").Value = "superadmin"
$ws.Range("I226").Value = "now()"

# Row 227
$ws.Range("A227").Formula = "=SUM(A226,1)"
$ws.Range("B227").Value = "preferredLang"
$ws.Range("C227").Value = "idioma preferido del usuario"
$ws.Range("D227").Value = "string"
$ws.Range("E227").Value = "{`"value`":`"தமிழ்`",`"code`":`"tam`"}"
$ws.Range("F227").Value = "spa"
$ws.Range("G227").Value = "TRUE"
$ws.Range("H227").Value = "superadmin"
$ws.Range("I227").Value = "now()"

# Row 228
$ws.Range("A228").Formula = "=SUM(A227,1)"
$ws.Range("B228").Value = "preferredLang"
$ws.Range("C228").Value = "ಬಳಕೆದಾರ ಆದ್ಯತೆಯ ಭಾಷೆ"
$ws.Range("D228").Value = "string"
$ws.Range("E228").Value = "{`"value`":`"Española`",`"code`":`"spa`"}"
$ws.Range("F228").Value = "kan"
$ws.Range("G228").Value = "TRUE"
$ws.Range("H228").Value = "superadmin"
$ws.Range("I228").Value = "now()"

# Row 229
$ws.Range("A229").Formula = "=SUM(A228,1)"
$ws.Range("B229").Value = "preferredLang"
$ws.Range("C229").Value = "ಬಳಕೆದಾರ ಆದ್ಯತೆಯ ಭಾಷೆ"
$ws.Range("D229").Value = "string"
$ws.Range("E229").Value = "{`"value`":`"ಕನ್ನಡ`",`"code`":`"kan`"}"
$ws.Range("F229").Value = "kan"
$ws.Range("G229").Value = "TRUE"
$ws.Range("H229").Value = "superadmin"
$ws.Range("I229").Value = "now()"

# Row 230
$ws.Range("A230").Formula = "=SUM(A229,1)"
$ws.Range("B230").Value = "preferredLang"
$ws.Range("C230").Value = "ಬಳಕೆದಾರ ಆದ್ಯತೆಯ ಭಾಷೆ"
$ws.Range("D230").Value = "string"
$ws.Range("E230").Value = "{`"value`":`"हिंदी`",`"code`":`"hin`"}"
$ws.Range("F230").Value = "kan"
$ws.Range("G230").Value = "TRUE"
$ws.Range("H230").Value = "superadmin"
$ws.Range("I230").Value = "now()"

# Row 231
$ws.Range("A231").Formula = "=SUM(A230,1)"
$ws.Range("B231").Value = "preferredLang"
$ws.Range("C231").Value = "ಬಳಕೆದಾರ ಆದ್ಯತೆಯ ಭಾಷೆ"
$ws.Range("D231").Value = "string"
$ws.Range("E231").Value = "{`"value`":`"தமிழ்`",`"code`":`"tam`"}"
$ws.Range("F231").Value = "kan"
$ws.Range("G231").Value = "TRUE"
$ws.Range("H231").Value = "superadmin"
$ws.Range("I231").Value = "now()"

# Row 232
$ws.Range("A232").Formula = "=SUM(A231,1)"
$ws.Range("B232").Value = "preferredLang"
$ws.Range("C232").Value = "उपयोगकर्ता पसंदीदा भाषा"
$ws.Range("D232").Value = "string"
$ws.Range("E232").Value = "{`"value`":`"Española`",`"code`":`"spa`"}"
$ws.Range("F232").Value = "hin"
$ws.Range("G232").Value = "TRUE"
$ws.Range("H232").Value = "superadmin"
$ws.Range("I232").Value = "now()"

# Row 233
$ws.Range("A233").Formula = "=SUM(A232,1)"
$ws.Range("B233").Value = "preferredLang"
$ws.Range("C233").Value = "उपयोगकर्ता पसंदीदा भाषा"
$ws.Range("D233").Value = "string"
$ws.Range("E233").Value = "{`"value`":`"ಕನ್ನಡ`",`"code`":`"kan`"}"
$ws.Range("F233").Value = "hin"
$ws.Range("G233").Value = "TRUE"
$ws.Range("H233").Value = "superadmin"
$ws.Range("I233").Value = "now()"

# Row 234
$ws.Range("A234").Formula = "=SUM(A233,1)"
$ws.Range("B234").Value = "preferredLang"
$ws.Range("C234").Value = "उपयोगकर्ता पसंदीदा भाषा"
$ws.Range("D234").Value = "string"
$ws.Range("E234").Value = "{`"value`":`"हिंदी`",`"code`":`"hin`"}"
$ws.Range("F234").Value = "hin"
$ws.Range("G234").Value = "TRUE"
$ws.Range("H234").Value = "superadmin"
$ws.Range("I234").Value = "now()"

# Row 235
$ws.Range("A235").Formula = "=SUM(A234,1)"
$ws.Range("B235").Value = "preferredLang"
$ws.Range("C235").Value = "उपयोगकर्ता पसंदीदा भाषा"
$ws.Range("D235").Value = "string"
$ws.Range("E235").Value = "{`"value`":`"தமிழ்`",`"code`":`"tam`"}"
$ws.Range("F235").Value = "hin"
$ws.Range("G235").Value = "TRUE"
$ws.Range("H235").Value = "superadmin"
$ws.Range("I235").Value = "now()"

# Row 236
$ws.Range("A236").Formula = "=SUM(A235,1)"
$ws.Range("B236").Value = "preferredLang"
$ws.Range("C236").Value = "பயனர் விருப்பமான மொழி"
$ws.Range("D236").Value = "string"
$ws.Range("E236").Value = "{`"value`":`"Española`",`"code`":`"spa`"}"
$ws.Range("F236").Value = "tam"
$ws.Range("G236").Value = "TRUE"
$ws.Range("H236").Value = "superadmin"
$ws.Range("I236").Value = "now()"

# Row 237
$ws.Range("A237").Formula = "=SUM(A236,1)"
$ws.Range("B237").Value = "preferredLang"
$ws.Range("C237").Value = "பயனர் விருப்பமான மொழி"
$ws.Range("D237").Value = "string"
$ws.Range("E237").Value = "{`"value`":`"ಕನ್ನಡ`",`"code`":`"kan`"}"
$ws.Range("F237").Value = "tam"
$ws.Range("G237").Value = "TRUE"
$ws.Range("H237").Value = "superadmin"
$ws.Range("I237").Value = "now()"

# Row 238
$ws.Range("A238").Formula = "=SUM(A237,1)"
$ws.Range("B238").Value = "preferredLang"
$ws.Range("C238").Value = "பயனர் விருப்பமான மொழி"
$ws.Range("D238").Value = "string"
$ws.Range("E238").Value = "{`"value`":`"हिंदी`",`"code`":`"hin`"}"
$ws.Range("F238").Value = "tam"
$ws.Range("G238").Value = "TRUE"
$ws.Range("H238").Value = "superadmin"
$ws.Range("I238").Value = "now()"

# Row 239
$ws.Range("A239").Formula = "=SUM(A238,1)"
$ws.Range("B239").Value = "preferredLang"
$ws.Range("C239").Value = "பயனர் விருப்பமான மொழி"
$ws.Range("D239").Value = "string"
$ws.Range("E239").Value = "{`"value`":`"தமிழ்`",`"code`":`"tam`"}"
$ws.Range("F239").Value = "tam"
$ws.Range("G239").Value = "TRUE"
$ws.Range("H239").Value = "superadmin"
$ws.Range("I239").Value = "now()"
